$wb = $excel.ActiveWorkbook

# Remove the "TestResultExcelFilePath" column (column H) from the two
# NI Scenario sheets that held it as an Input-sheet output-path column.
$ws3 = $wb.Worksheets.Item("ProcessPayrollForNIMonthly")
$ws3.Columns.Item(8).EntireColumn.Delete()

$ws4 = $wb.Worksheets.Item("TestReports")
$ws4.Columns.Item(8).EntireColumn.Delete()
$ws4.Range("J8").Select()

$ws3.Activate()
$ws3.Range("H1:H1048576").Select()
